$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1446.0303
$ws.Range("I40").Value = 1683.25
$ws.Range("J40").Value = 1413.3103
$ws.Range("K40").Value = 1683.25
$ws.Range("L40").Value = 1413.3103
$ws.Range("M40").Value = -1508.25
$ws.Range("N40").Value = -1763.3103

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 735.7143
$ws.Range("I103").Value = 975
$ws.Range("K103").Value = 2925
$ws.Range("M103").Value = -2339

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2704175
$ws.Range("I137").Value = 5556571
$ws.Range("J137").Value = 1904.8948
$ws.Range("K137").Value = 16669713
$ws.Range("L137").Value = 5714.6844
$ws.Range("M137").Value = -16667163
$ws.Range("N137").Value = -10814.6844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 25693458
$ws.Range("I61").Value = 29442194
$ws.Range("J61").Value = 202056
$ws.Range("K61").Value = 29442194
$ws.Range("L61").Value = 202056
$ws.Range("M61").Value = -29441982
$ws.Range("N61").Value = -202480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7000990.5
$ws.Range("I74").Value = 7172447.5
$ws.Range("J74").Value = 1000000
$ws.Range("K74").Value = 7172447.5
$ws.Range("L74").Value = 1000000
$ws.Range("M74").Value = -7171573.5
$ws.Range("N74").Value = -1001748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 37987.5
$ws.Range("J76").Value = 37987.5
$ws.Range("L76").Value = 37987.5
$ws.Range("N76").Value = -38663.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7000990.5
$ws.Range("I77").Value = 7172447.5
$ws.Range("J77").Value = 1000000
$ws.Range("K77").Value = 35862237.5
$ws.Range("L77").Value = 5000000
$ws.Range("M77").Value = -35857869.5
$ws.Range("N77").Value = -5008736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 37987.5
$ws.Range("J79").Value = 37987.5
$ws.Range("L79").Value = 37987.5
$ws.Range("N79").Value = -40327.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 79453.3
$ws.Range("I132").Value = 54013.74
$ws.Range("K132").Value = 162041.22
$ws.Range("M132").Value = -159511.22

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 25693458
$ws.Range("I136").Value = 29442194
$ws.Range("J136").Value = 202056
$ws.Range("K136").Value = 88326582
$ws.Range("L136").Value = 606168
$ws.Range("M136").Value = -88324032
$ws.Range("N136").Value = -611268

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 39271.332
$ws.Range("J76").Value = 39271.332
$ws.Range("L76").Value = 39271.332
$ws.Range("N76").Value = -39901.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 39271.332
$ws.Range("J79").Value = 39271.332
$ws.Range("L79").Value = 39271.332
$ws.Range("N79").Value = -41455.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2716.2104
$ws.Range("I134").Value = 2356.077
$ws.Range("J134").Value = 3496.5
$ws.Range("K134").Value = 7068.231000000001
$ws.Range("L134").Value = 10489.5
$ws.Range("M134").Value = -4533.231000000001
$ws.Range("N134").Value = -15559.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2046.6129
$ws.Range("I31").Value = 952.1905
$ws.Range("J31").Value = 4344.9
$ws.Range("K31").Value = 952.1905
$ws.Range("L31").Value = 4344.9
$ws.Range("M31").Value = -657.1905
$ws.Range("N31").Value = -4934.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2046.6129
$ws.Range("I34").Value = 952.1905
$ws.Range("J34").Value = 4344.9
$ws.Range("K34").Value = 952.1905
$ws.Range("L34").Value = 4344.9
$ws.Range("M34").Value = -750.1905
$ws.Range("N34").Value = -4748.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17858998
$ws.Range("I58").Value = 23257264
$ws.Range("J58").Value = 3191.4614
$ws.Range("K58").Value = 23257264
$ws.Range("L58").Value = 3191.4614
$ws.Range("M58").Value = -23257061
$ws.Range("N58").Value = -3597.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 29625.111
$ws.Range("I132").Value = 1746.037
$ws.Range("J132").Value = 113262.336
$ws.Range("K132").Value = 5238.111
$ws.Range("L132").Value = 339787.008
$ws.Range("M132").Value = -2708.111
$ws.Range("N132").Value = -344847.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 33585.176
$ws.Range("I134").Value = 1977
$ws.Range("J134").Value = 136311.75
$ws.Range("K134").Value = 5931
$ws.Range("L134").Value = 408935.25
$ws.Range("M134").Value = -3396
$ws.Range("N134").Value = -414005.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 17858998
$ws.Range("I136").Value = 23257264
$ws.Range("J136").Value = 3191.4614
$ws.Range("K136").Value = 69771792
$ws.Range("L136").Value = 9574.3842
$ws.Range("M136").Value = -69769242
$ws.Range("N136").Value = -14674.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 71453.84
$ws.Range("J101").Value = 71453.84
$ws.Range("L101").Value = 71453.84
$ws.Range("N101").Value = -77943.84

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 144317.72
$ws.Range("I132").Value = 100865
$ws.Range("J132").Value = 252949.5
$ws.Range("K132").Value = 302595
$ws.Range("L132").Value = 758848.5
$ws.Range("M132").Value = -300065
$ws.Range("N132").Value = -763908.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3263.75
$ws.Range("I16").Value = 1101.75
$ws.Range("J16").Value = 9749.75
$ws.Range("K16").Value = 1101.75
$ws.Range("L16").Value = 9749.75
$ws.Range("M16").Value = -931.75
$ws.Range("N16").Value = -10089.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 29963.818
$ws.Range("J103").Value = 29963.818
$ws.Range("L103").Value = 29963.818
$ws.Range("N103").Value = -32307.818

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18655.45
$ws.Range("I132").Value = 1160.081
$ws.Range("J132").Value = 49480.617
$ws.Range("K132").Value = 3480.242999999999
$ws.Range("L132").Value = 148441.851
$ws.Range("M132").Value = -950.2429999999995
$ws.Range("N132").Value = -153501.851

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 75264.664
$ws.Range("I136").Value = 53660.105
$ws.Range("J136").Value = 126575.5
$ws.Range("K136").Value = 160980.315
$ws.Range("L136").Value = 379726.5
$ws.Range("M136").Value = -158430.315
$ws.Range("N136").Value = -384826.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 55468.973
$ws.Range("I132").Value = 40731.76
$ws.Range("J132").Value = 86171.5
$ws.Range("K132").Value = 122195.28
$ws.Range("L132").Value = 258514.5
$ws.Range("M132").Value = -119665.28
$ws.Range("N132").Value = -263574.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 40378.63
$ws.Range("I136").Value = 22706.738
$ws.Range("K136").Value = 68120.21400000001
$ws.Range("M136").Value = -65570.21400000001
